# add NIIT - simplify quick start
#
# Both the "John" and "Sally" sheets track a year-by-year projection
# starting at the current year. The example is refreshed so the
# projection now starts in 2020 instead of 2025: five new rows
# (2020-2024) are inserted right after the header row on each sheet,
# pushing the existing 2025+ rows down by five rows. The newly
# inserted rows pick up the same currency-style formatting used by
# the rest of the table, and the previously selected ranges/tabs are
# refreshed to reflect where each sheet's author was last working.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "John"
$ws2 = $wb.Worksheets.Item(2)   # "Sally"

function Prepend-Years {
    param($ws)

    # Push everything down 5 rows, starting right after the header.
    $ws.Rows("2:6").Insert()

    # The inserted rows default to the header's bold style; re-apply the
    # plain currency formatting used throughout the rest of the table by
    # copying formats from an existing, already-correct data row.
    $ws.Range("A8:I8").Copy() | Out-Null
    $ws.Range("A2:I6").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    # Fill in the five new years.
    $ws.Cells.Item(2, 1).Value = 2020
    $ws.Cells.Item(3, 1).Value = 2021
    $ws.Cells.Item(4, 1).Value = 2022
    $ws.Cells.Item(5, 1).Value = 2023
    $ws.Cells.Item(6, 1).Value = 2024
}

Prepend-Years $ws1
Prepend-Years $ws2

# Refresh each sheet's selection to match where the author left off.
$ws1.Range("A2:XFD6").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B11").Select() | Out-Null
